# edit.ps1 -- applies the three-part change described by the diff:
#   1. Remove the old "_GoBack" bookmark that sat after
#      ", lo añadimos al index.js".
#   2. Split "...incluido el remitente." into three runs, appending
#      ", lo añadimos en el index.js" before the trailing period (which
#      becomes its own run).
#   3. Split "Y en el lado del cliente, cuando capturemos un chat message"
#      into four runs, inserting " (index.html)" after "cliente" and
#      re-adding the "_GoBack" bookmark right after " (index.html" --
#      while keeping the (unchanged) trailing " " + "evento, ..." runs
#      of that paragraph intact as their own runs.

$d = $word.ActiveDocument
$nbsp = [char]0xA0

# ---------------------------------------------------------------------
# Change 1: drop the bookmark that originally followed
# ", lo añadimos al index.js"
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Change 2: "...incluido el remitente." ->
#   run1: "...incluido el remitente"
#   run2: ", lo añadimos en el index.js"
#   run3: "."
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("incluido el remitente.", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "incluido el remitente", 2)

$find2 = $d.Content
$find2.Find.Execute("incluido el remitente")
$insertPoint = $d.Range($find2.End, $find2.End)
$insertPoint.InsertAfter(", lo añadimos en el index.js.")

# Re-find the newly inserted text (without the trailing period) and
# toggle Bold on/off -- this forces the engine to keep it as its own
# run (with unchanged run properties) instead of merging it back into
# its identically-formatted neighbours.
$find3 = $d.Content
$find3.Find.Execute(", lo añadimos en el index.js")
$find3.Bold = 1
$find3.Bold = 0

# ---------------------------------------------------------------------
# Change 3: "Y en el lado del cliente, cuando capturemos un chat message" ->
#   run1: "Y en el lado del cliente"
#   run2: " (index.html"
#   [bookmarkStart/End _GoBack]
#   run3: ")"
#   run4: ", cuando capturemos un chat message"
#   (run5: " " and run6: "evento, ..." stay as in the original)
# ---------------------------------------------------------------------

# Locate the original run precisely (it contains a non-breaking space
# between "un" and "chat").
$findA = $d.Content
$findA.Find.Execute("Y en el lado del cliente, cuando capturemos")
$runAStart = $findA.Start
$origRunAText = "Y en el lado del cliente, cuando capturemos un" + $nbsp + "chat message"
$runAEnd = $runAStart + $origRunAText.Length

# Insert " (index.html)" right after "...cliente".
$insertPos = $runAStart + "Y en el lado del cliente".Length
$insertRange = $d.Range($insertPos, $insertPos)
$insertText = " (index.html)"
$insertRange.InsertAfter($insertText)
$shift = $insertText.Length

# The (unchanged) single-space run and the "evento..." run that used to
# follow run A now sit right after the inserted text; locate them so we
# can re-isolate them (inserting text merges same-formatted neighbours
# into a single run).
$runBStart = $runAEnd + $shift
$runBEnd = $runBStart + 1
$runCStart = $runBEnd
$origRunCText = "evento, lo incluiremos en la p" + [char]0xE1 + "gina." + $nbsp + `
                "El c" + [char]0xF3 + "digo total de JavaScript del lado del cliente ahora equivale a:"
$runCEnd = $runCStart + $origRunCText.Length

# Split off run C ("evento...") first.
$rC = $d.Range($runCStart, $runCEnd)
$rC.Bold = 1
$rC.Bold = 0

# Split off run B (the lone space) next.
$rB = $d.Range($runBStart, $runBEnd)
$rB.Bold = 1
$rB.Bold = 0

# Boundaries of the four new sub-runs inside [runAStart, runBStart).
$p1End = $insertPos                                  # end of "Y en el lado del cliente"
$p2End = $insertPos + " (index.html".Length           # end of " (index.html"
$p3Start = $p2End
$p3End = $p3Start + 1                                 # ")"
$p4Start = $p3End
$p4End = $runBStart                                   # ", cuando capturemos un<nbsp>chat message"

# Split rightmost-first so earlier boundaries remain valid.
$r4 = $d.Range($p4Start, $p4End)
$r4.Bold = 1
$r4.Bold = 0

$r3 = $d.Range($p3Start, $p3End)
$r3.Bold = 1
$r3.Bold = 0

$r2 = $d.Range($p1End, $p2End)
$r2.Bold = 1
$r2.Bold = 0

# Re-add the "_GoBack" bookmark between run2 (" (index.html") and run3 (")").
$bmRange = $d.Range($p2End, $p2End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
